$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet: Property1 -> DataNode
$ws.Name = "DataNode"

# Column width tweaks (B,C,D,E,G,H,I,J got ~1/16-pt narrower in the source file).
# ColumnWidth is quantized by the host to 1/7-character steps, so we pick the
# COM value whose round-trip lands closest to the authored width.
$ws.Columns.Item(2).ColumnWidth = 10.857142857142858   # -> 11.625 target
$ws.Columns.Item(3).ColumnWidth = 16.428571428571427   # -> 17.125 target
$ws.Columns.Item(4).ColumnWidth = 21.857142857142858   # -> 22.625 target
$ws.Columns.Item(5).ColumnWidth = 16.428571428571427   # -> 17.125 target
$ws.Columns.Item(7).ColumnWidth = 11.857142857142858   # -> 12.625 target
$ws.Columns.Item(8).ColumnWidth = 17.714285714285715   # -> 18.375 target
$ws.Columns.Item(9).ColumnWidth = 18.714285714285715   # -> 19.375 target
$ws.Columns.Item(10).ColumnWidth = 13.142857142857142  # -> 13.875 target

# Row height tweaks
$ws.Rows.Item(1).RowHeight = 27
$ws.Rows.Item(8).RowHeight = 81

# Selection moved from A9 to H13 (still in the frozen bottom-left pane)
$ws.Range("H13").Select()
